$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell T1 = "17" (week 17), kept as text like the other header cells
$ws.Range("T1").NumberFormat = "@"
$ws.Range("T1").Value = "17"
$ws.Range("T1").Font.Bold = $true
$ws.Range("T1").HorizontalAlignment = -4108
$ws.Range("T1").NumberFormat = "General"

# Correct an existing data value: Q26 changes from 2 to 0
$ws.Range("Q26").Value = 0

# New week-17 counts in column T for each data row
$ws.Range("T2").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("T6").Value = 21
$ws.Range("T7").Value = 3
$ws.Range("T8").Value = 29
$ws.Range("T9").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("T12").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("T16").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("T19").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("T21").Value = 0
$ws.Range("T22").Value = 0
$ws.Range("T23").Value = 4
$ws.Range("T24").Value = 0
$ws.Range("T25").Value = 0
$ws.Range("T26").Value = 0
$ws.Range("T27").Value = 6
$ws.Range("T28").Value = 39
$ws.Range("T29").Value = 0
$ws.Range("T30").Value = 0
$ws.Range("T32").Value = 52
$ws.Range("T33").Value = 1
$ws.Range("T34").Value = 0
$ws.Range("T35").Value = 0
$ws.Range("T36").Value = 0
$ws.Range("T37").Value = 0
$ws.Range("T38").Value = 0
$ws.Range("T39").Value = 0
$ws.Range("T41").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("T43").Value = 0
$ws.Range("T44").Value = 0
$ws.Range("T45").Value = 0
$ws.Range("T46").Value = 0
$ws.Range("T47").Value = 0
$ws.Range("T48").Value = 0
$ws.Range("T49").Value = 1
$ws.Range("T50").Value = 0
$ws.Range("T51").Value = 0
$ws.Range("T52").Value = 0
$ws.Range("T53").Value = 0
